$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Update Tür (type) for D203 and D204 rows from "Derslik" to "Lab"
$ws.Range("B4").Value = "Lab"
$ws.Range("B5").Value = "Lab"

# Update Kapasite (capacity) values
$ws.Range("C4").Value = 80
$ws.Range("C5").Value = 90
$ws.Range("C6").Value = 70
$ws.Range("C7").Value = 70
$ws.Range("C8").Value = 80
$ws.Range("C9").Value = 80
$ws.Range("C10").Value = 60
$ws.Range("C11").Value = 50

# Move active selection to C6, matching final saved view state
$ws.Range("C6").Select()
